$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.356.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4682'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2851'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06543'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.83'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07934'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '98.01'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.877.09'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.153'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6831'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '282.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.351.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.97'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.446'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.121.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.54%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007322'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.171'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.73'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.188'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.949'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.391'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09752'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.431'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.484'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.126'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04738'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.139'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7144'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.725'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01874'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.361'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.557'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.984'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8546'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4214'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.249'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '960.57'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.344'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1132'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.88%  '
